$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Price column as Text so numeric-looking values (e.g. "0.999", "1.70")
# are preserved verbatim as text instead of being re-interpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "39.742.25"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "2.213.33"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "292.30"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").Value = "86.57"
$ws.Range("E6").Value = "  +6.85%  "
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "0.471"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").Value = "30.19"
$ws.Range("E10").Value = "  +3.91%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.0785"
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("D12").Value = "47.35"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").Value = "0.109"
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").Value = "6.33"
$ws.Range("E14").Value = "  +1.90%  "
$ws.Range("D15").Value = "2.554.58"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "14.03"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").Value = "2.204.68"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "0.727"
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("D19").Value = "39.685.61"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("D20").Value = "11.53"
$ws.Range("E20").Value = "  +12.96%  "
$ws.Range("D21").Value = "0.0₃0879"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").Value = "5.80"
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("D23").Value = "65.82"
$ws.Range("E23").Value = "  +2.09%  "
$ws.Range("D24").Value = "235.47"
$ws.Range("E24").Value = "  +4.02%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  +2.78%  "
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("D28").Value = "22.71"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("D30").Value = "9.26"
$ws.Range("E30").Value = "  +2.58%  "
$ws.Range("D31").Value = "32.73"
$ws.Range("E31").Value = "  +4.02%  "
$ws.Range("D32").Value = "151.98"
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "4.93"
$ws.Range("E34").Value = "  +3.14%  "
$ws.Range("D35").Value = "0.0717"
$ws.Range("E35").Value = "  +3.51%  "
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "2.79"
$ws.Range("E38").Value = "  +6.64%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "15.98"
$ws.Range("E39").Value = "  +4.61%  "
$ws.Range("D40").Value = "0.0985"
$ws.Range("E40").Value = "  +2.88%  "
$ws.Range("D41").Value = "1.70"
$ws.Range("E41").Value = "  +3.16%  "
$ws.Range("D42").Value = "2.081.01"
$ws.Range("E42").Value = "  +9.68%  "
$ws.Range("D43").Value = "3.78"
$ws.Range("E43").Value = "  +5.13%  "
$ws.Range("E44").Value = "  +6.02%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0269"
$ws.Range("E45").Value = "  +3.90%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "9.99"
$ws.Range("E46").Value = "  +11.20%  "
$ws.Range("D47").Value = "17.65"
$ws.Range("E47").Value = "  +10.07%  "
$ws.Range("D48").Value = "2.62"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").Value = "2.426.18"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("D50").Value = "70.75"
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("D51").Value = "89.03"
$ws.Range("E51").Value = "  +2.10%  "
